$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDe = "Wir werden Ihnen nun einige Reihen von jeweils vier einzelnen Wörtern oder kurzen Ausdrücken präsentieren, z.B. Teller – Schachtel – Ratte – Parkett. Bitte lesen Sie diese Reihen aufmerksam und achten Sie dabei besonders auf den Sprachrhythmus, d.h. auf die Betonungsmuster der Wörter. <br/> In manchen dieser Reihen weisen alle Elemente, d.h. alle enthaltenen Einzelwörter oder Mehrwort-Ausdrücke, dasselbe Betonungsmuster auf.  In einigen dieser Reihen ist jedoch genau ein Element enthalten, das ein anderes Betonungsmuster aufweist als die übrigen – im Beispiel oben ist es das  Wort „Parkett“. Ihre Aufgabe besteht darin, für jede Reihe zu beurteilen, ob ein Element mit abweichendem Betonungsmuster enthalten ist oder nicht. Falls alle Elemente der Reihe das gleiche Betonungsmuster aufweisen, klicken Sie bitte auf das Feld „alle gleich“. Falls ein Element ein abweichendes Betonungsmuster aufweist, klicken Sie bitte auf das Feld dieses abweichenden Elements.<br/> Bitte beachten Sie, dass Sie für jede Reihe 10 Sekunden Zeit haben, um Ihr Urteil abzugeben."

$newEn = "We will now present you with some series of four individual words or short expressions each, e.g. plate - box - rat - parquet.  Please read these series carefully, paying particular  attention to the rhythm of speech, i.e. the stress patterns of the words. <br/> In some of these series, however, there is exactly one element that has a different stress pattern than the others - in the example above it is the word `"parquet`".  Your task is to judge for each row whether or not there is an element with a different stress pattern.  If all the elements in the row have the same stress pattern, please click on the field `"all the same`".  If an element has a different stress pattern, please click on the field of this different element. <br/>Please note that you have 10 seconds for each row to make your judgement."

$ws.Range("B3").Value = $newDe
$ws.Range("C3").Value = $newEn

$ws.Range("C4").Select()
